# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# Both sheets list the same events (the latter aggregates every category),
# so the same old -> new value pairs apply, just at different row offsets.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value, for the "展览" sheet
$exhibitionUpdates = @{
    2  = 15176
    3  = 19554
    5  = 167
    13 = 63
    14 = 214
    17 = 1520
    20 = 113
    22 = 8207
    23 = 993
    24 = 44
    27 = 1272
    28 = 19
    31 = 6538
    35 = 158
    37 = 5590
    38 = 1018
    39 = 29
    41 = 63
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Row -> new value, for the "全部类型" sheet
$allTypesUpdates = @{
    2  = 15176
    3  = 19554
    5  = 167
    13 = 63
    14 = 214
    17 = 1520
    21 = 113
    23 = 8207
    24 = 993
    25 = 44
    28 = 1272
    29 = 19
    34 = 6538
    38 = 158
    40 = 5590
    41 = 1018
    42 = 29
    44 = 63
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
